$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SMP")
$ws2 = $wb.Worksheets.Item("Lobby")

# Add the new "Velocity" sheet after "Lobby", mirroring the header row
# (Nume / Link / Tip / Github) used on the other two sheets.
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Velocity"
$ws3.Range("A1").Value = "Nume"
$ws3.Range("B1").Value = "Link"
$ws3.Range("C1").Value = "Tip"
$ws3.Range("D1").Value = "Github"
$ws3.Range("A2").Select() | Out-Null

# Turn the existing link text in column B (rows 3 and 6) of the SMP sheet
# into real hyperlinks, pointing at the URL already shown in each cell
# (this also creates the "Hyperlink" cell style used by the linked cells).
$ws1.Hyperlinks.Add($ws1.Range("B3"), $ws1.Range("B3").Text) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B6"), $ws1.Range("B6").Text) | Out-Null

# Update the saved selections: Lobby now just has A1:D1 selected (and is no
# longer the active tab), while SMP becomes the active tab with B6 selected.
$ws2.Range("A1:D1").Select() | Out-Null
$ws1.Range("B6").Select() | Out-Null

# Make SMP the active sheet/tab of the workbook.
$ws1.Activate() | Out-Null

Write-Output "done"
